$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the source data ---
# Row 26 = "RM 232" is deleted outright.
$ws.Rows(26).Delete()
# After that delete, the old row 28 ("SC 92") is now row 27.
$ws.Rows(27).Delete()

# --- Column F (column 6) value swaps for the rows above the deleted block ---
$ws.Cells.Item(6, 6).Value = 16.43
$ws.Cells.Item(8, 6).Value = ""
$ws.Cells.Item(18, 6).Value = 18.35
$ws.Cells.Item(20, 6).Value = ""
$ws.Cells.Item(23, 6).Value = 16.48
$ws.Cells.Item(25, 6).Value = ""

# --- Cell-level corrections within the rows that shifted up after the deletions ---
# Row 27 = "SC 101": column D (4) now has a value
$ws.Cells.Item(27, 4).Value = -14.6
# Row 28 = "SC 105": column D now blank
$ws.Cells.Item(28, 4).Value = ""
# Row 29 = "SC 119": column D now blank
$ws.Cells.Item(29, 4).Value = ""
# Row 30 = "SC 120": column D now has a value, column F now has a value
$ws.Cells.Item(30, 4).Value = -13.6
$ws.Cells.Item(30, 6).Value = 16.89
# Row 32 = "SC 193": column D now blank
$ws.Cells.Item(32, 4).Value = ""
